{"js": "const tables = context.document.body.tables;\ntables.load(\"items\");\nawait context.sync();\nconst table = tables.items[0];\n\n// Each lattice-multiplication cell holds 5 lines joined by vertical-tab\n// (Word line-break) characters: the problem, the two digit factors,\n// the separator, and the two partial-product placeholders.\nconst cellValues = [\n  [0, 0, \"76 x 55\\v  5    5\\v  ----\\v7|    |\\v6|    |\"],\n  [0, 1, \"13 x 77\\v  7    7\\v  ----\\v1|    |\\v3|    |\"],\n  [0, 2, \"63 x 10\\v  1    0\\v  ----\\v6|    |\\v3|    |\"],\n  [1, 0, \"50 x 94\\v  9    4\\v  ----\\v5|    |\\v0|    |\"],\n  [1, 1, \"42 x 28\\v  2    8\\v  ----\\v4|    |\\v2|    |\"],\n  [1, 2, \"31 x 17\\v  1    7\\v  ----\\v3|    |\\v1|    |\"],\n  [2, 0, \"30 x 44\\v  4    4\\v  ----\\v3|    |\\v0|    |\"],\n  [2, 1, \"80 x 41\\v  4    1\\v  ----\\v8|    |\\v0|    |\"],\n  [2, 2, \"16 x 52\\v  5    2\\v  ----\\v1|    |\\v6|    |\"],\n  [3, 0, \"40 x 84\\v  8    4\\v  ----\\v4|    |\\v0|    |\"],\n  [3, 1, \"21 x 25\\v  2    5\\v  ----\\v2|    |\\v1|    |\"],\n  [3, 2, \"12 x 20\\v  2    0\\v  ----\\v1|    |\\v2|    |\"],\n  [4, 0, \"28 x 96\\v  9    6\\v  ----\\v2|    |\\v8|    |\"],\n  [4, 1, \"11 x 92\\v  9    2\\v  ----\\v1|    |\\v1|    |\"],\n  [4, 2, \"50 x 89\\v  8    9\\v  ----\\v5|    |\\v0|    |\"],\n];\n\nfor (const [row, col, text] of cellValues) {\n  const cell = table.getCell(row, col);\n  const paragraphs = cell.body.paragraphs;\n  paragraphs.load(\"items\");\n  await context.sync();\n  // Replace the text of the cell's (single) paragraph in place so the\n  // existing run (and its 32-half-point font size) is preserved.\n  paragraphs.items[0].insertText(text, Word.InsertLocation.replace);\n  await context.sync();\n}\n", "ps1": "$d = $word.ActiveDocument\n$t = $d.Tables.Item(1)\n\n# Each lattice-multiplication cell holds 5 lines joined by a vertical-tab\n# (Word line-break) character: the problem, the two digit factors, the\n# separator, and the two partial-product placeholders.\n# Rows/columns below are 1-based, matching Word's COM Cell(row, col).\n$t.Cell(1, 1).Range.Text = \"76 x 55`v  5    5`v  ----`v7|    |`v6|    |\"\n$t.Cell(1, 2).Range.Text = \"13 x 77`v  7    7`v  ----`v1|    |`v3|    |\"\n$t.Cell(1, 3).Range.Text = \"63 x 10`v  1    0`v  ----`v6|    |`v3|    |\"\n$t.Cell(2, 1).Range.Text = \"50 x 94`v  9    4`v  ----`v5|    |`v0|    |\"\n$t.Cell(2, 2).Range.Text = \"42 x 28`v  2    8`v  ----`v4|    |`v2|    |\"\n$t.Cell(2, 3).Range.Text = \"31 x 17`v  1    7`v  ----`v3|    |`v1|    |\"\n$t.Cell(3, 1).Range.Text = \"30 x 44`v  4    4`v  ----`v3|    |`v0|    |\"\n$t.Cell(3, 2).Range.Text = \"80 x 41`v  4    1`v  ----`v8|    |`v0|    |\"\n$t.Cell(3, 3).Range.Text = \"16 x 52`v  5    2`v  ----`v1|    |`v6|    |\"\n$t.Cell(4, 1).Range.Text = \"40 x 84`v  8    4`v  ----`v4|    |`v0|    |\"\n$t.Cell(4, 2).Range.Text = \"21 x 25`v  2    5`v  ----`v2|    |`v1|    |\"\n$t.Cell(4, 3).Range.Text = \"12 x 20`v  2    0`v  ----`v1|    |`v2|    |\"\n$t.Cell(5, 1).Range.Text = \"28 x 96`v  9    6`v  ----`v2|    |`v8|    |\"\n$t.Cell(5, 2).Range.Text = \"11 x 92`v  9    2`v  ----`v1|    |`v1|    |\"\n$t.Cell(5, 3).Range.Text = \"50 x 89`v  8    9`v  ----`v5|    |`v0|    |\"\n"}
